$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 26 new patient-id rows (rows 16-41) below the existing data (A1:A15).
# One of the new rows (A20) mirrors the pre-existing blank PatientID row (A2),
# so it is entered the same way Excel represents a manually-typed empty/text
# cell - a leading apostrophe - which keeps it sharing the existing empty
# shared string instead of minting a new one; the quote-prefix style picked
# up along the way is reset back to Normal so the cell stays plain.

$ws.Cells.Item(16, 1).Value = "PEP_ID-2005243"
$ws.Cells.Item(17, 1).Value = "PEP_ID-2005248"
$ws.Cells.Item(18, 1).Value = "PEP_ID-2005250"
$ws.Cells.Item(19, 1).Value = "PEP_ID-2005275"

$ws.Cells.Item(20, 1).Value = "'"
$ws.Cells.Item(20, 1).Style = "Normal"

$ws.Cells.Item(21, 1).Value = "PEP_ID-2005396"
$ws.Cells.Item(22, 1).Value = "null"
$ws.Cells.Item(23, 1).Value = "PEP_ID-2005410"
$ws.Cells.Item(24, 1).Value = "PEP_ID-2005419"
$ws.Cells.Item(25, 1).Value = "PEP_ID-2005424"
$ws.Cells.Item(26, 1).Value = "PEP_ID-2005427"
$ws.Cells.Item(27, 1).Value = "null"
$ws.Cells.Item(28, 1).Value = "null"
$ws.Cells.Item(29, 1).Value = "PEP_ID-2005478"
$ws.Cells.Item(30, 1).Value = "PEP_ID-2005485"
$ws.Cells.Item(31, 1).Value = "PEP_ID-2005493"
$ws.Cells.Item(32, 1).Value = "PEP_ID-2005559"
$ws.Cells.Item(33, 1).Value = "PEP_ID-2005562"
$ws.Cells.Item(34, 1).Value = "PEP_ID-2005568"
$ws.Cells.Item(35, 1).Value = "PEP_ID-2005585"
$ws.Cells.Item(36, 1).Value = "PEP_ID-2005588"
$ws.Cells.Item(37, 1).Value = "PEP_ID-2005591"
$ws.Cells.Item(38, 1).Value = "PEP_ID-2005596"
$ws.Cells.Item(39, 1).Value = "PEP_ID-2005603"
$ws.Cells.Item(40, 1).Value = "PEP_ID-2005607"
$ws.Cells.Item(41, 1).Value = "PEP_ID-2005609"
